$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3220
$ws.Range("I64").Value = 3350
$ws.Range("J64").Value = 3133.3333
$ws.Range("K64").Value = 3350
$ws.Range("L64").Value = 3133.3333
$ws.Range("M64").Value = -3102
$ws.Range("N64").Value = -3629.3333

$ws.Range("H67").Value = 3220
$ws.Range("I67").Value = 3350
$ws.Range("J67").Value = 3133.3333
$ws.Range("K67").Value = 3350
$ws.Range("L67").Value = 3133.3333
$ws.Range("M67").Value = -2492
$ws.Range("N67").Value = -4849.3333

$ws.Range("H98").Value = 1567.1177
$ws.Range("I98").Value = 1743.2858
$ws.Range("J98").Value = 745
$ws.Range("K98").Value = 1743.2858
$ws.Range("L98").Value = 745
$ws.Range("M98").Value = -245.2858000000001
$ws.Range("N98").Value = -3741

$ws.Range("H122").Value = 1567.1177
$ws.Range("I122").Value = 1743.2858
$ws.Range("J122").Value = 745
$ws.Range("K122").Value = 5229.857400000001
$ws.Range("L122").Value = 2235
$ws.Range("M122").Value = -2779.857400000001
$ws.Range("N122").Value = -7135

$ws.Range("H135").Value = 481.1875
$ws.Range("I135").Value = 481.1875
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 4330.6875
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 409.66666
$ws.Range("I5").Value = 191.6
$ws.Range("K5").Value = 191.6
$ws.Range("M5").Value = -79.59999999999999

$ws.Range("H32").Value = 4444.4287
$ws.Range("I32").Value = 2904.8
$ws.Range("J32").Value = 11287.223
$ws.Range("K32").Value = 2904.8
$ws.Range("L32").Value = 11287.223
$ws.Range("M32").Value = -2617.8
$ws.Range("N32").Value = -11861.223

$ws.Range("H45").Value = 1548.625
$ws.Range("I45").Value = 1198.25
$ws.Range("K45").Value = 1198.25
$ws.Range("M45").Value = -821.25

$ws.Range("H74").Value = 2239.5264
$ws.Range("I74").Value = 1735.8182
$ws.Range("J74").Value = 2932.125
$ws.Range("K74").Value = 1735.8182
$ws.Range("L74").Value = 2932.125
$ws.Range("M74").Value = -861.8181999999999
$ws.Range("N74").Value = -4680.125

$ws.Range("H77").Value = 2239.5264
$ws.Range("I77").Value = 1735.8182
$ws.Range("J77").Value = 2932.125
$ws.Range("K77").Value = 8679.091
$ws.Range("L77").Value = 14660.625
$ws.Range("M77").Value = -4311.091
$ws.Range("N77").Value = -23396.625

$ws.Range("H122").Value = 1153.4
$ws.Range("I122").Value = 1148.2222
$ws.Range("K122").Value = 3444.6666
$ws.Range("M122").Value = -994.6665999999996

$ws.Range("H132").Value = 1891.4688
$ws.Range("I132").Value = 1218.826
$ws.Range("K132").Value = 3656.478
$ws.Range("M132").Value = -1126.478

$ws.Range("H133").Value = 32753.334
$ws.Range("J133").Value = 32753.334
$ws.Range("L133").Value = 32753.334
$ws.Range("N133").Value = -37813.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 409.66666
$ws.Range("I4").Value = 191.6
$ws.Range("K4").Value = 191.6
$ws.Range("M4").Value = -76.59999999999999

$ws.Range("H131").Value = 55780
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 55780
$ws.Range("K131").Value = 0
$ws.Range("N131").Value = -65860
$ws.Range("M131").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2531.1667
$ws.Range("I62").Value = 2531.1667
$ws.Range("K62").Value = 2531.1667
$ws.Range("M62").Value = -1907.1667

$ws.Range("H65").Value = 2531.1667
$ws.Range("I65").Value = 2531.1667
$ws.Range("K65").Value = 12655.8335
$ws.Range("M65").Value = -9535.833500000001

$ws.Range("H99").Value = 2666.5557
$ws.Range("I99").Value = 1850
$ws.Range("J99").Value = 2899.8572
$ws.Range("K99").Value = 1850
$ws.Range("L99").Value = 2899.8572
$ws.Range("M99").Value = -352
$ws.Range("N99").Value = -5895.8572

$ws.Range("H122").Value = 2076
$ws.Range("I122").Value = 942
$ws.Range("J122").Value = 10014
$ws.Range("K122").Value = 2826
$ws.Range("L122").Value = 30042
$ws.Range("M122").Value = -376
$ws.Range("N122").Value = -34942

$ws.Range("H126").Value = 2666.5557
$ws.Range("I126").Value = 1850
$ws.Range("J126").Value = 2899.8572
$ws.Range("K126").Value = 5550
$ws.Range("L126").Value = 8699.571599999999
$ws.Range("M126").Value = -3080
$ws.Range("N126").Value = -13639.5716

$ws.Range("H134").Value = 1105.5581
$ws.Range("I134").Value = 1078.75
$ws.Range("K134").Value = 3236.25
$ws.Range("M134").Value = -701.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 61.333332
$ws.Range("I6").Value = 61.333332
$ws.Range("K6").Value = 183.999996
$ws.Range("M6").Value = -70.99999600000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4502
$ws.Range("I70").Value = 4702.8
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 4702.8
$ws.Range("L70").Value = 4000
$ws.Range("M70").Value = -4432.8
$ws.Range("N70").Value = -4540

$ws.Range("H73").Value = 4502
$ws.Range("I73").Value = 4702.8
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 4702.8
$ws.Range("L73").Value = 4000
$ws.Range("M73").Value = -3766.8
$ws.Range("N73").Value = -5872

$ws.Range("H102").Value = 2450.3157
$ws.Range("I102").Value = 2475.389
$ws.Range("J102").Value = 1999
$ws.Range("K102").Value = 2475.389
$ws.Range("L102").Value = 1999
$ws.Range("M102").Value = -853.3890000000001
$ws.Range("N102").Value = -5243

$ws.Range("H122").Value = 2372.182
$ws.Range("I122").Value = 1849.25
$ws.Range("J122").Value = 2671
$ws.Range("K122").Value = 5547.75
$ws.Range("L122").Value = 8013
$ws.Range("M122").Value = -3097.75
$ws.Range("N122").Value = -12913

$ws.Range("H132").Value = 3207621.2
$ws.Range("I132").Value = 5496365.5
$ws.Range("K132").Value = 16489096.5
$ws.Range("M132").Value = -16486566.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3596.5
$ws.Range("I7").Value = 2316.4
$ws.Range("J7").Value = 9997
$ws.Range("K7").Value = 2316.4
$ws.Range("L7").Value = 9997
$ws.Range("M7").Value = -2204.4
$ws.Range("N7").Value = -10221

$ws.Range("H22").Value = 2619.2307
$ws.Range("J22").Value = 1787.5
$ws.Range("L22").Value = 1787.5
$ws.Range("N22").Value = -2377.5

$ws.Range("H27").Value = 2619.2307
$ws.Range("J27").Value = 1787.5
$ws.Range("L27").Value = 1787.5
$ws.Range("N27").Value = -2001.5

$ws.Range("H46").Value = 2782.3333
$ws.Range("I46").Value = 1933.3334
$ws.Range("J46").Value = 3631.3333
$ws.Range("K46").Value = 1933.3334
$ws.Range("L46").Value = 3631.3333
$ws.Range("M46").Value = -1745.3334
$ws.Range("N46").Value = -4007.3333

$ws.Range("H126").Value = 3596.5
$ws.Range("I126").Value = 2316.4
$ws.Range("J126").Value = 9997
$ws.Range("K126").Value = 6949.200000000001
$ws.Range("L126").Value = 29991
$ws.Range("M126").Value = -4479.200000000001
$ws.Range("N126").Value = -34931

$ws.Range("H136").Value = 3855
$ws.Range("I136").Value = 2843.5715
$ws.Range("K136").Value = 8530.7145
$ws.Range("M136").Value = -5980.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1649.75
$ws.Range("I81").Value = 1366.3334
$ws.Range("K81").Value = 2732.6668
$ws.Range("M81").Value = -1671.6668

$ws.Range("H84").Value = 1649.75
$ws.Range("I84").Value = 1366.3334
$ws.Range("K84").Value = 13663.334
$ws.Range("M84").Value = -8359.333999999999

$ws.Range("H100").Value = 806.4286
$ws.Range("I100").Value = 461.25
$ws.Range("K100").Value = 922.5
$ws.Range("M100").Value = -381.5

$ws.Range("H122").Value = 49654.625
$ws.Range("I122").Value = 49654.625
$ws.Range("K122").Value = 148963.875
$ws.Range("M122").Value = -146513.875

$ws.Range("H126").Value = 4559
$ws.Range("J126").Value = 11250
$ws.Range("L126").Value = 33750
$ws.Range("N126").Value = -38690

$ws.Range("H132").Value = 1871.84
$ws.Range("I132").Value = 1248.1177
$ws.Range("J132").Value = 3197.25
$ws.Range("K132").Value = 3744.3531
$ws.Range("L132").Value = 9591.75
$ws.Range("M132").Value = -1214.3531
$ws.Range("N132").Value = -14651.75
